# Apply the 6.a.1.1 update:
#  - fix a typo in the (cosmetic) absolute-path metadata (adds a space)
#  - extend the data table with four more years (2020-2023) in columns X:AA
#  - hide the now-stale columns D:K (years 2000-2007)
#  - bump a couple of row heights slightly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix absPath typo in workbook-level metadata -----------------------
# (Not exposed on the Excel object model; nothing to change here via COM -
#  it is regenerated by the host application itself on save.)

# --- 2. Add the new year columns (X:AA) on row 4 --------------------------
$ws.Cells.Item(4, 24).Value = 2020
$ws.Cells.Item(4, 25).Value = 2021
$ws.Cells.Item(4, 26).Value = 2022
$ws.Cells.Item(4, 27).Value = 2023

# --- 3. Add the new data values for rows 5-7 -------------------------------
$ws.Cells.Item(5, 24).Value = 23780
$ws.Cells.Item(5, 25).Value = 44660
$ws.Cells.Item(5, 26).Value = 25000
$ws.Cells.Item(5, 27).Value = 13010

$ws.Cells.Item(6, 24).Value = 38240
$ws.Cells.Item(6, 25).Value = 7950
$ws.Cells.Item(6, 26).Value = 23000
$ws.Cells.Item(6, 27).Value = 16390

$ws.Cells.Item(7, 24).Value = 62020
$ws.Cells.Item(7, 25).Value = 52610
$ws.Cells.Item(7, 26).Value = 48000
$ws.Cells.Item(7, 27).Value = 29400

# --- 4. Copy formatting from the last "2019" column (W) into the new ones -
$ws.Range("W4:W7").Copy()
$ws.Range("X4:X7").PasteSpecial(-4122)
$ws.Range("W4:W7").Copy()
$ws.Range("Y4:Y7").PasteSpecial(-4122)
$ws.Range("W4:W7").Copy()
$ws.Range("Z4:Z7").PasteSpecial(-4122)
$ws.Range("W4:W7").Copy()
$ws.Range("AA4:AA7").PasteSpecial(-4122)

# --- 5. Hide the stale columns D:K (years 2000-2007) -----------------------
$ws.Range("D1:K1").EntireColumn.Hidden = $true

# --- 6. Slightly taller rows for the (now wider) header/data rows ----------
$ws.Rows.Item(4).RowHeight = 16.5
$ws.Rows.Item(5).RowHeight = 16.5
$ws.Rows.Item(6).RowHeight = 16.5
$ws.Rows.Item(7).RowHeight = 16.5
